# Edit script: apply "Ran code for averaged intensities on spiral schemes"
# - Adds 3 new shared strings for spiral sampling schemes
# - Re-runs averaged-intensity computation, changing row 10-16 values and header refs
# - Adds 3 new rows (17-19) for the HexGrid schemes that got pushed down

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row 2 (C2:P2): shared string indices shift due to 3 new strings inserted ---
$ws.Range("C2").Value = "[1, 1, 1]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 2, 0]"
$ws.Range("F2").Value = "[3, 1, 1]"
$ws.Range("G2").Value = "[2, 2, 2]"
$ws.Range("H2").Value = "[4, 0, 0]"
$ws.Range("I2").Value = "[3, 3, 1]"
$ws.Range("J2").Value = "[4, 2, 0]"
$ws.Range("K2").Value = "[4, 2, 2]"
$ws.Range("L2").Value = "[5, 1, 1]"
$ws.Range("M2").Value = "[3, 3, 3]"
$ws.Range("N2").Value = "2Pairs"
$ws.Range("O2").Value = "4Pairs"
$ws.Range("P2").Value = "MaxUnique"

# --- Update existing data rows 10-16 (re-ordered / recomputed values) ---
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.560627062854345
$ws.Range("D10").Value = 1.039049773750829
$ws.Range("E10").Value = 1.622924357132956
$ws.Range("F10").Value = 0.7297912232260775
$ws.Range("G10").Value = 1.560627062854345
$ws.Range("H10").Value = 1.039049773750829
$ws.Range("I10").Value = 1.015516809015217
$ws.Range("J10").Value = 1.111178102939896
$ws.Range("K10").Value = 0.8657733098031337
$ws.Range("L10").Value = 0.7045161732785121
$ws.Range("M10").Value = 1.560627062854345
$ws.Range("N10").Value = 1.330987065441892
$ws.Range("O10").Value = 1.238098104241052
$ws.Range("P10").Value = 1.081172101500121

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.3266449798664107
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 2.707299763610842
$ws.Range("F11").Value = 0.6322686641608687
$ws.Range("G11").Value = 0.3266449798664107
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 1.845512433833404
$ws.Range("J11").Value = 0.659302194500726
$ws.Range("K11").Value = 1.194430008683945
$ws.Range("L11").Value = 0.1079776849009707
$ws.Range("M11").Value = 0.3266449798664107
$ws.Range("N11").Value = 1.353649881805421
$ws.Range("O11").Value = 0.9165533519095305
$ws.Range("P11").Value = 0.9341794661946459

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.326945119996509
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 2.697373161718997
$ws.Range("F12").Value = 0.6326850678874875
$ws.Range("G12").Value = 0.326945119996509
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 1.846151887686994
$ws.Range("J12").Value = 0.659350189919743
$ws.Range("K12").Value = 1.194828812756972
$ws.Range("L12").Value = 0.1077581603305528
$ws.Range("M12").Value = 0.326945119996509
$ws.Range("N12").Value = 1.348686580859498
$ws.Range("O12").Value = 0.9142508374007483
$ws.Range("P12").Value = 0.9331365500371569

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.3261271920021528
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 2.708262603382333
$ws.Range("F13").Value = 0.6311318209268607
$ws.Range("G13").Value = 0.3261271920021528
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 1.847009888419409
$ws.Range("J13").Value = 0.6598362696975579
$ws.Range("K13").Value = 1.19182743123876
$ws.Range("L13").Value = 0.1074267990116627
$ws.Range("M13").Value = 0.3261271920021528
$ws.Range("N13").Value = 1.354131301691167
$ws.Range("O13").Value = 0.9163804040778367
$ws.Range("P13").Value = 0.9339527505848422

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.6152760000000006
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 1.291092000000003
$ws.Range("F14").Value = 1.216532000000001
$ws.Range("G14").Value = 0.6152760000000006
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 1.123308
$ws.Range("J14").Value = 0.374812
$ws.Range("K14").Value = 2.400911999999999
$ws.Range("L14").Value = 0.2579080000000003
$ws.Range("M14").Value = 0.6152760000000006
$ws.Range("N14").Value = 0.6455460000000014
$ws.Range("O14").Value = 0.780725000000001
$ws.Range("P14").Value = 0.9099800000000005

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.88
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0.02
$ws.Range("F15").Value = 1.750999999999997
$ws.Range("G15").Value = 0.88
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0.4598625
$ws.Range("J15").Value = 0.12
$ws.Range("K15").Value = 3.517224999999991
$ws.Range("L15").Value = 0.4002999999999995
$ws.Range("M15").Value = 0.88
$ws.Range("N15").Value = 0.01
$ws.Range("O15").Value = 0.6627499999999993
$ws.Range("P15").Value = 0.8935484374999986

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.9692009281535994
$ws.Range("D16").Value = 0.3768308365312038
$ws.Range("E16").Value = 0.4238952747008025
$ws.Range("F16").Value = 1.421300689510394
$ws.Range("G16").Value = 0.9692009281535994
$ws.Range("H16").Value = 0.3768308365312038
$ws.Range("I16").Value = 0.6957360648192029
$ws.Range("J16").Value = 0.5095861814272005
$ws.Range("K16").Value = 2.387224659046406
$ws.Range("L16").Value = 0.6529274897408005
$ws.Range("M16").Value = 0.9691247108095998
$ws.Range("N16").Value = 0.4003630556160032
$ws.Range("O16").Value = 0.7978069322239998
$ws.Range("P16").Value = 0.9295877654912011

# --- Add new rows 17-19 (previously rows 13 & 15, and row 14, pushed down) ---
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9770890042618984
$ws.Range("D17").Value = 0.9958155504106795
$ws.Range("E17").Value = 1.001773448469824
$ws.Range("F17").Value = 0.9850964846689166
$ws.Range("G17").Value = 0.9770890042618984
$ws.Range("H17").Value = 0.9958155504106795
$ws.Range("I17").Value = 0.9926289592092368
$ws.Range("J17").Value = 0.9984010437941702
$ws.Range("K17").Value = 0.990076030960036
$ws.Range("L17").Value = 0.9897808431246279
$ws.Range("M17").Value = 0.977044084176299
$ws.Range("N17").Value = 0.9987944994402516
$ws.Range("O17").Value = 0.9899436219528295
$ws.Range("P17").Value = 0.9913326706124236
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.271676545846232
$ws.Range("D18").Value = 1.128058125257921
$ws.Range("E18").Value = 0.9383261268166483
$ws.Range("F18").Value = 1.171945002309921
$ws.Range("G18").Value = 1.271676545846232
$ws.Range("H18").Value = 1.128058125257921
$ws.Range("I18").Value = 0.9126687736956245
$ws.Range("J18").Value = 0.8557232636279815
$ws.Range("K18").Value = 0.9909354238778219
$ws.Range("L18").Value = 1.09571072153018
$ws.Range("M18").Value = 1.271709783152289
$ws.Range("N18").Value = 1.033192126037285
$ws.Range("O18").Value = 1.127501450057681
$ws.Range("P18").Value = 1.045630497870291
$ws.Range("A16").Copy()
$ws.Range("A18").PasteSpecial(-4122)

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9043493467296733
$ws.Range("D19").Value = 0.9928018094940526
$ws.Range("E19").Value = 0.9281815662111584
$ws.Range("F19").Value = 1.046100518163948
$ws.Range("G19").Value = 0.9043493467296733
$ws.Range("H19").Value = 0.9928018094940526
$ws.Range("I19").Value = 0.9137071422386869
$ws.Range("J19").Value = 1.051581539425303
$ws.Range("K19").Value = 0.9770193670041363
$ws.Range("L19").Value = 1.076258636503809
$ws.Range("M19").Value = 0.9042591566025288
$ws.Range("N19").Value = 0.9604916878526055
$ws.Range("O19").Value = 0.9678583101497081
$ws.Range("P19").Value = 0.986249990721346
$ws.Range("A16").Copy()
$ws.Range("A19").PasteSpecial(-4122)

$excel.CutCopyMode = 0